# Updates the cryptos price/volume(1h) table to the latest scraped
# snapshot (GitHub Actions run on Fri May 17 01:44:11 UTC 2024).
#
# The Price column (D) holds values formatted/scraped as plain text
# (e.g. "65.388.46", "161.70", "0.0₃0996") rather than numbers, so a
# handful of cells look number-like (e.g. "0.999", "161.70") and would
# otherwise be auto-coerced into real numbers by Excel's Value setter
# (losing the trailing zero / exact text). For those we force text entry
# with a leading apostrophe and then ClearFormats() to drop the implicit
# "quote prefix" text number-format Excel applies, keeping the cell style
# identical to the untouched cells around it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '65.388.46'
$ws.Cells.Item(2, 5).Value = '  -0.86%  '
$ws.Cells.Item(3, 4).Value = '2.948.58'
$ws.Cells.Item(3, 5).Value = '  -2.26%  '
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.10%  '
$ws.Cells.Item(5, 4).Value = "'570.32"
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -1.99%  '
$ws.Cells.Item(6, 4).Value = "'161.70"
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +1.03%  '
$ws.Cells.Item(7, 4).Value = "'0.999"
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -0.16%  '
$ws.Cells.Item(8, 5).Value = '  -0.36%  '
$ws.Cells.Item(9, 4).Value = '2.947.37'
$ws.Cells.Item(9, 5).Value = '  -2.20%  '
$ws.Cells.Item(10, 5).Value = '  -4.21%  '
$ws.Cells.Item(11, 4).Value = "'0.151"
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -2.38%  '
$ws.Cells.Item(12, 4).Value = "'0.458"
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +0.40%  '
$ws.Cells.Item(13, 4).Value = "'0.0000244"
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  -2.28%  '
$ws.Cells.Item(14, 4).Value = "'34.58"
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  -0.40%  '
$ws.Cells.Item(15, 5).Value = '  -1.16%  '
$ws.Cells.Item(16, 4).Value = '65.323.23'
$ws.Cells.Item(16, 5).Value = '  -1.10%  '
$ws.Cells.Item(17, 4).Value = '3.435.71'
$ws.Cells.Item(17, 5).Value = '  -2.35%  '
$ws.Cells.Item(18, 5).Value = '  +0.92%  '
$ws.Cells.Item(19, 4).Value = '2.944.32'
$ws.Cells.Item(19, 5).Value = '  -2.41%  '
$ws.Cells.Item(20, 4).Value = "'15.83"
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +14.32%  '
$ws.Cells.Item(21, 4).Value = "'444.52"
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -3.09%  '
$ws.Cells.Item(22, 4).Value = "'0.698"
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +1.64%  '
$ws.Cells.Item(23, 4).Value = "'7.29"
$ws.Cells.Item(23, 4).ClearFormats()
$ws.Cells.Item(23, 5).Value = '  -1.04%  '
$ws.Cells.Item(24, 4).Value = "'82.58"
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +0.39%  '
$ws.Cells.Item(25, 4).Value = "'2.25"
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -0.46%  '
$ws.Cells.Item(26, 4).Value = "'12.14"
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  -2.29%  '
$ws.Cells.Item(27, 2).Value = 'Dai'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(27, 4).Value = "'1.00"
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.05%  '
$ws.Cells.Item(28, 2).Value = 'RenderToken'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(28, 4).Value = "'9.98"
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  -5.96%  '
$ws.Cells.Item(29, 4).Value = "'2.46"
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +4.18%  '
$ws.Cells.Item(30, 4).Value = "'8.06"
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -1.52%  '
$ws.Cells.Item(31, 5).Value = '  -0.87%  '
$ws.Cells.Item(32, 4).Value = '0.0₃0996'
$ws.Cells.Item(32, 5).Value = '  -5.03%  '
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(33, 4).Value = "'0.113"
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +2.09%  '
$ws.Cells.Item(34, 2).Value = 'EthereumClassic'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(34, 4).Value = "'27.24"
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +0.98%  '
$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -0.18%  '
$ws.Cells.Item(36, 5).Value = '  -2.54%  '
$ws.Cells.Item(37, 4).Value = "'5.74"
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -0.65%  '
$ws.Cells.Item(38, 4).Value = "'49.48"
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  -1.40%  '
$ws.Cells.Item(39, 4).Value = "'45.22"
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +3.56%  '
$ws.Cells.Item(40, 5).Value = '  -0.32%  '
$ws.Cells.Item(41, 4).Value = "'0.121"
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  -0.67%  '
$ws.Cells.Item(42, 4).Value = "'2.82"
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -5.54%  '
$ws.Cells.Item(43, 4).Value = "'1.93"
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -10.61%  '
$ws.Cells.Item(44, 4).Value = "'8.54"
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +0.58%  '
$ws.Cells.Item(45, 4).Value = "'385.25"
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +0.78%  '
$ws.Cells.Item(46, 4).Value = "'0.0351"
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -1.38%  '
$ws.Cells.Item(47, 4).Value = '2.683.93'
$ws.Cells.Item(47, 5).Value = '  -3.82%  '
$ws.Cells.Item(48, 5).Value = '  -0.01%  '
$ws.Cells.Item(49, 5).Value = '  +0.05%  '
$ws.Cells.Item(50, 4).Value = "'2.18"
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  +2.92%  '
$ws.Cells.Item(51, 4).Value = "'23.55"
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -1.15%  '
